$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of an existing data row (row 19, which carries the
# exact style set needed: name/type/expr1 cells boxed, expr2-3 + action plain)
# onto the new row 29, then fill in the values for the new
# "DocumentRepository - Default assignee" rule.
$ws.Range("B19:H19").Copy()
$ws.Range("B29:H29").PasteSpecial(-4122)
$ws.Rows.Item(29).RowHeight = 30

# Fill in the new-string cells in the same order the shared-string table
# records them (Rule Name, Expression 2, Action, then Expression 1), so the
# appended <si> entries line up with the authored workbook.
$ws.Range("B29").Value = "DocumentRepository " + [char]8211 + " Default assignee"
$ws.Range("E29").Value = "repositoryType == 'PERSONAL'"
$ws.Range("H29").Value = "assignee, creator"
$ws.Range("D29").Value = "creator != null"
$ws.Range("C29").Value = "DOC_REPO"
$ws.Range("F29").Value = ""
$ws.Range("G29").Value = ""

# Update the active sheet view / selection to match the authored state
$excel.ActiveWindow.ScrollRow = 16
[void]$ws.Range("D29").Select()

$excel.Iteration = $true
$excel.MaxChange = 0.0001
